# Auto-commit data refresh: update last-charge-end timestamps (col D) for all rows,
# and replace the terminal alert snapshot (cols A-D) for rows 18-48, including three
# newly appended rows (46-48). Selection moves to F10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRefreshTime = 45990.276759259257

# Rows 2-17: only the "last charge end time" (column D) refreshes.
$ws.Range("D2:D17").Value = $newRefreshTime

# Rows 18-48: full snapshot rows (station, terminal, alert-since time, refresh time).
$rows = @(
    @{ Row = 18; A = "长沙特来电飞狐四方坪南区充电站"; B = "101号直流"; C = 45987.552604166667; D = 45990.276759259257 }
    @{ Row = 19; A = "长沙特来电飞狐四方坪西区充电站"; B = "904号直流"; C = 45987.576736111114; D = 45990.276759259257 }
    @{ Row = 20; A = "长沙特来电飞狐四方坪南区充电站"; B = "103号直流"; C = 45988.086469907408; D = 45990.276759259257 }
    @{ Row = 21; A = "长沙特来电飞狐四方坪西区充电站"; B = "505号直流"; C = 45988.235497685186; D = 45990.276759259257 }
    @{ Row = 22; A = "长沙特来电飞狐四方坪西区充电站"; B = "B01号直流"; C = 45988.522523148145; D = 45990.276759259257 }
    @{ Row = 23; A = "长沙特来电飞狐四方坪西区充电站"; B = "B02号直流"; C = 45989.073136574072; D = 45990.276759259257 }
    @{ Row = 24; A = "长沙特来电飞狐四方坪西区充电站"; B = "802号直流"; C = 45989.194456018522; D = 45990.276759259257 }
    @{ Row = 25; A = "长沙特来电飞狐四方坪东区充电站"; B = "402号直流"; C = 45989.208715277775; D = 45990.276759259257 }
    @{ Row = 26; A = "长沙特来电飞狐四方坪西区充电站"; B = "502号直流"; C = 45989.209733796299; D = 45990.276759259257 }
    @{ Row = 27; A = "长沙市开福区高岭香江国际城充电站建设项目"; B = "111号直流"; C = 45989.310312499998; D = 45990.276759259257 }
    @{ Row = 28; A = "长沙市开福区高岭香江国际城充电站建设项目"; B = "108号直流"; C = 45989.402465277781; D = 45990.276759259257 }
    @{ Row = 29; A = "长沙特来电飞狐四方坪南区充电站"; B = "105号直流"; C = 45989.537094907406; D = 45990.276759259257 }
    @{ Row = 30; A = "长沙市开福区高岭香江国际城充电站建设项目"; B = "109号直流"; C = 45989.540393518517; D = 45990.276759259257 }
    @{ Row = 31; A = "长沙特来电飞狐四方坪西区充电站"; B = "402号直流"; C = 45989.54146990741; D = 45990.276759259257 }
    @{ Row = 32; A = "长沙特来电飞狐四方坪西区充电站"; B = "503号直流"; C = 45989.544745370367; D = 45990.276759259257 }
    @{ Row = 33; A = "长沙特来电飞狐四方坪西区充电站"; B = "603号直流"; C = 45989.545405092591; D = 45990.276759259257 }
    @{ Row = 34; A = "长沙特来电飞狐四方坪南区充电站"; B = "201号直流"; C = 45989.545717592591; D = 45990.276759259257 }
    @{ Row = 35; A = "长沙特来电飞狐四方坪西区充电站"; B = "702号直流"; C = 45989.545902777776; D = 45990.276759259257 }
    @{ Row = 36; A = "长沙市开福区高岭香江国际城充电站建设项目"; B = "102号直流"; C = 45989.557164351849; D = 45990.276759259257 }
    @{ Row = 37; A = "长沙特来电飞狐四方坪西区充电站"; B = "903号直流"; C = 45989.55777777778; D = 45990.276759259257 }
    @{ Row = 38; A = "长沙市开福区高岭香江国际城充电站建设项目"; B = "112号直流"; C = 45989.558182870373; D = 45990.276759259257 }
    @{ Row = 39; A = "长沙特来电飞狐四方坪东区充电站"; B = "103号直流"; C = 45989.570868055554; D = 45990.276759259257 }
    @{ Row = 40; A = "长沙特来电飞狐四方坪东区充电站"; B = "006A号直流"; C = 45989.582245370373; D = 45990.276759259257 }
    @{ Row = 41; A = "长沙市开福区高岭香江国际城充电站建设项目"; B = "305号直流"; C = 45989.587905092594; D = 45990.276759259257 }
    @{ Row = 42; A = "长沙市开福区高岭香江国际城充电站建设项目"; B = "208号直流"; C = 45989.590081018519; D = 45990.276759259257 }
    @{ Row = 43; A = "长沙特来电飞狐四方坪南区充电站"; B = "405号直流"; C = 45989.595104166663; D = 45990.276759259257 }
    @{ Row = 44; A = "长沙特来电飞狐四方坪南区充电站"; B = "401号直流"; C = 45989.607418981483; D = 45990.276759259257 }
    @{ Row = 45; A = "长沙市开福区高岭香江国际城充电站建设项目"; B = "206号直流"; C = 45989.661643518521; D = 45990.276759259257 }
    @{ Row = 46; A = "长沙特来电飞狐四方坪东区充电站"; B = "001B号直流"; C = 45989.666354166664; D = 45990.276759259257 }
    @{ Row = 47; A = "长沙特来电飞狐四方坪西区充电站"; B = "401号直流"; C = 45989.704826388886; D = 45990.276759259257 }
    @{ Row = 48; A = "长沙特来电飞狐四方坪西区充电站"; B = "604号直流"; C = 45989.766516203701; D = 45990.276759259257 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
}

# Move the active selection to F10 to match the saved view state.
$ws.Range("F10").Select()
